$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 970.6667
$ws.Range("J17").Value = 970.6667
$ws.Range("L17").Value = 2912.0001
$ws.Range("N17").Value = -3248.0001
$ws.Range("H33").Value = 219.17647
$ws.Range("I33").Value = 123.07692
$ws.Range("K33").Value = 123.07692
$ws.Range("M33").Value = 105.92308
$ws.Range("H40").Value = 2471.2856
$ws.Range("J40").Value = 1914
$ws.Range("L40").Value = 1914
$ws.Range("N40").Value = -2264
$ws.Range("H51").Value = 1713.9445
$ws.Range("I51").Value = 950.3333
$ws.Range("J51").Value = 2095.75
$ws.Range("K51").Value = 950.3333
$ws.Range("L51").Value = 2095.75
$ws.Range("M51").Value = -466.3333
$ws.Range("N51").Value = -3063.75
$ws.Range("H86").Value = 4408.5654
$ws.Range("I86").Value = 3929.8
$ws.Range("K86").Value = 3929.8
$ws.Range("M86").Value = -2806.8
$ws.Range("H89").Value = 4408.5654
$ws.Range("I89").Value = 3929.8
$ws.Range("K89").Value = 19649
$ws.Range("M89").Value = -14033
$ws.Range("H107").Value = 1720.9642
$ws.Range("I107").Value = 1623
$ws.Range("K107").Value = 1623
$ws.Range("M107").Value = 297
$ws.Range("H132").Value = 5055184.5
$ws.Range("I132").Value = 6538756
$ws.Range("J132").Value = 11040.733
$ws.Range("K132").Value = 19616268
$ws.Range("L132").Value = 33122.199
$ws.Range("M132").Value = -19613738
$ws.Range("N132").Value = -38182.199
$ws.Range("H137").Value = 1101.2329
$ws.Range("I137").Value = 901.3077
$ws.Range("J137").Value = 1330.5588
$ws.Range("K137").Value = 2703.9231
$ws.Range("L137").Value = 3991.6764
$ws.Range("M137").Value = -153.9231
$ws.Range("N137").Value = -9091.6764
$ws.Range("H138").Value = 613728.7
$ws.Range("I138").Value = 961.6896400000001
$ws.Range("J138").Value = 1036829.7
$ws.Range("K138").Value = 2885.06892
$ws.Range("L138").Value = 3110489.1
$ws.Range("M138").Value = 2254.93108
$ws.Range("N138").Value = -3120769.1
$ws.Range("H141").Value = 550.4474
$ws.Range("I141").Value = 550.4474
$ws.Range("K141").Value = 1651.3422
$ws.Range("M141").Value = 3528.6578

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 21262.8
$ws.Range("I2").Value = 1503.6666
$ws.Range("J2").Value = 50901.5
$ws.Range("K2").Value = 1503.6666
$ws.Range("L2").Value = 50901.5
$ws.Range("M2").Value = -1390.6666
$ws.Range("N2").Value = -51127.5
$ws.Range("H32").Value = 4544.4375
$ws.Range("I32").Value = 4158.4287
$ws.Range("K32").Value = 4158.4287
$ws.Range("M32").Value = -3871.4287
$ws.Range("H61").Value = 17242248
$ws.Range("I61").Value = 18868662
$ws.Range("K61").Value = 18868662
$ws.Range("M61").Value = -18868450
$ws.Range("H74").Value = 1147
$ws.Range("I74").Value = 789.4194
$ws.Range("J74").Value = 2378.6667
$ws.Range("K74").Value = 789.4194
$ws.Range("L74").Value = 2378.6667
$ws.Range("M74").Value = 84.5806
$ws.Range("N74").Value = -4126.6667
$ws.Range("H77").Value = 1147
$ws.Range("I77").Value = 789.4194
$ws.Range("J77").Value = 2378.6667
$ws.Range("K77").Value = 3947.097
$ws.Range("L77").Value = 11893.3335
$ws.Range("M77").Value = 420.9030000000002
$ws.Range("N77").Value = -20629.3335
$ws.Range("H116").Value = 21262.8
$ws.Range("I116").Value = 1503.6666
$ws.Range("J116").Value = 50901.5
$ws.Range("K116").Value = 1503.6666
$ws.Range("L116").Value = 50901.5
$ws.Range("M116").Value = 790.3334
$ws.Range("N116").Value = -55489.5
$ws.Range("H131").Value = 49463.332
$ws.Range("J131").Value = 49463.332
$ws.Range("L131").Value = 49463.332
$ws.Range("N131").Value = -59543.332
$ws.Range("H132").Value = 2948.2173
$ws.Range("I132").Value = 3263.0625
$ws.Range("J132").Value = 2228.5715
$ws.Range("K132").Value = 9789.1875
$ws.Range("L132").Value = 6685.7145
$ws.Range("M132").Value = -7259.1875
$ws.Range("N132").Value = -11745.7145
$ws.Range("H136").Value = 17242248
$ws.Range("I136").Value = 18868662
$ws.Range("K136").Value = 56605986
$ws.Range("M136").Value = -56603436

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 21262.8
$ws.Range("I3").Value = 1503.6666
$ws.Range("J3").Value = 50901.5
$ws.Range("K3").Value = 1503.6666
$ws.Range("L3").Value = 50901.5
$ws.Range("M3").Value = -1389.6666
$ws.Range("N3").Value = -51129.5
$ws.Range("H107").Value = 2090.0908
$ws.Range("I107").Value = 1647.375
$ws.Range("J107").Value = 3270.6667
$ws.Range("K107").Value = 1647.375
$ws.Range("L107").Value = 3270.6667
$ws.Range("M107").Value = 272.625
$ws.Range("N107").Value = -7110.6667
$ws.Range("H134").Value = 8319.105
$ws.Range("I134").Value = 1820.2307
$ws.Range("J134").Value = 22400
$ws.Range("K134").Value = 5460.6921
$ws.Range("L134").Value = 67200
$ws.Range("M134").Value = -2925.6921
$ws.Range("N134").Value = -72270

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 71429830
$ws.Range("I16").Value = 83334550
$ws.Range("J16").Value = 1495
$ws.Range("K16").Value = 83334550
$ws.Range("L16").Value = 1495
$ws.Range("M16").Value = -83334263
$ws.Range("N16").Value = -2069
$ws.Range("H94").Value = 815.5
$ws.Range("I94").Value = 806
$ws.Range("J94").Value = 820.25
$ws.Range("K94").Value = 806
$ws.Range("L94").Value = 820.25
$ws.Range("M94").Value = -355
$ws.Range("N94").Value = -1722.25
$ws.Range("H113").Value = 71429830
$ws.Range("I113").Value = 83334550
$ws.Range("J113").Value = 1495
$ws.Range("K113").Value = 83334550
$ws.Range("L113").Value = 1495
$ws.Range("M113").Value = -83332380
$ws.Range("N113").Value = -5835
$ws.Range("H132").Value = 5830.125
$ws.Range("I132").Value = 5232.4287
$ws.Range("J132").Value = 10014
$ws.Range("K132").Value = 15697.2861
$ws.Range("L132").Value = 30042
$ws.Range("M132").Value = -13167.2861
$ws.Range("N132").Value = -35102
$ws.Range("H135").Value = 32740
$ws.Range("J135").Value = 32740
$ws.Range("L135").Value = 32740
$ws.Range("N135").Value = -42880

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 886.4706
$ws.Range("J68").Value = 833.3333
$ws.Range("L68").Value = 2499.9999
$ws.Range("N68").Value = -4121.9999
$ws.Range("H71").Value = 886.4706
$ws.Range("J71").Value = 833.3333
$ws.Range("L71").Value = 7499.9997
$ws.Range("N71").Value = -15611.9997
$ws.Range("H114").Value = 452.86365
$ws.Range("I114").Value = 343.8
$ws.Range("J114").Value = 543.75
$ws.Range("K114").Value = 1031.4
$ws.Range("L114").Value = 1631.25
$ws.Range("M114").Value = 2222.6
$ws.Range("N114").Value = -8139.25
$ws.Range("H131").Value = 18183092
$ws.Range("J131").Value = 1337
$ws.Range("L131").Value = 4011
$ws.Range("N131").Value = -14091
$ws.Range("H132").Value = 1398.8572
$ws.Range("I132").Value = 1396.75
$ws.Range("J132").Value = 1401.6666
$ws.Range("K132").Value = 12570.75
$ws.Range("L132").Value = 12614.9994
$ws.Range("M132").Value = -10040.75
$ws.Range("N132").Value = -17674.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4966.3335
$ws.Range("I122").Value = 4900
$ws.Range("J122").Value = 4999.5
$ws.Range("K122").Value = 14700
$ws.Range("L122").Value = 14998.5
$ws.Range("M122").Value = -12250
$ws.Range("N122").Value = -19898.5
$ws.Range("H130").Value = 37340
$ws.Range("J130").Value = 37340
$ws.Range("L130").Value = 37340
$ws.Range("N130").Value = -47380
$ws.Range("H132").Value = 1976.35
$ws.Range("I132").Value = 1862.6111
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 5587.8333
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -3057.8333
$ws.Range("N132").Value = -14060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 265.94873
$ws.Range("I55").Value = 211.35
$ws.Range("J55").Value = 323.42105
$ws.Range("K55").Value = 211.35
$ws.Range("L55").Value = 323.42105
$ws.Range("M55").Value = -38.34999999999999
$ws.Range("N55").Value = -669.4210499999999
$ws.Range("H132").Value = 18134.117
$ws.Range("I132").Value = 1127.919
$ws.Range("J132").Value = 45491.914
$ws.Range("K132").Value = 3383.757000000001
$ws.Range("L132").Value = 136475.742
$ws.Range("M132").Value = -853.7570000000005
$ws.Range("N132").Value = -141535.742

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 15000
$ws.Range("J70").Value = 15000
$ws.Range("L70").Value = 15000
$ws.Range("N70").Value = -15630
$ws.Range("H73").Value = 15000
$ws.Range("J73").Value = 15000
$ws.Range("L73").Value = 15000
$ws.Range("N73").Value = -17184
$ws.Range("H81").Value = 400
$ws.Range("I81").Value = 335.42856
$ws.Range("K81").Value = 670.85712
$ws.Range("M81").Value = 390.14288
$ws.Range("H84").Value = 400
$ws.Range("I84").Value = 335.42856
$ws.Range("K84").Value = 3354.2856
$ws.Range("M84").Value = 1949.7144
$ws.Range("H107").Value = 436.0625
$ws.Range("I107").Value = 392.07693
$ws.Range("J107").Value = 626.6667
$ws.Range("K107").Value = 1176.23079
$ws.Range("L107").Value = 1880.0001
$ws.Range("M107").Value = 743.7692099999999
$ws.Range("N107").Value = -5720.0001
$ws.Range("H122").Value = 83338210
$ws.Range("I122").Value = 96159280
$ws.Range("J122").Value = 1252.5
$ws.Range("K122").Value = 288477840
$ws.Range("L122").Value = 3757.5
$ws.Range("M122").Value = -288475390
$ws.Range("N122").Value = -8657.5
$ws.Range("H132").Value = 8635.777
$ws.Range("I132").Value = 12285.333
$ws.Range("J132").Value = 1336.6666
$ws.Range("K132").Value = 36855.999
$ws.Range("L132").Value = 4009.9998
$ws.Range("M132").Value = -34325.999
$ws.Range("N132").Value = -9069.9998
$ws.Range("H136").Value = 665.64514
$ws.Range("J136").Value = 1345.8572
$ws.Range("L136").Value = 4037.5716
$ws.Range("N136").Value = -9137.571599999999

